$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" (summary) sheet: insert a new row for the 2022-Q4 quarter right
#    after the header row, pushing all existing quarter rows down by one.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# Copy formatting only from the row just below (the original row2, now
# shifted to row3) so the new row matches the existing look (bold/bordered
# index column, plain data columns).
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0

# ---------------------------------------------------------------------------
# 2. Add the new "2022-Q4" detail sheet. Duplicate the structurally
#    identical "2022-Q3" sheet (same columns/styles) and place the copy
#    immediately before it, then overwrite with the Q4 fund holding data.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3, $null)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# The source sheet has two data rows; the Q4 sheet only needs one.
$q4.Rows.Item(3).Delete()

$q4.Range("A2").Value = 0

$q4.Range("B2").NumberFormat = "@"
$q4.Range("B2").Value = "562530"
$q4.Range("B2").ClearFormats()

$q4.Range("C2").NumberFormat = "@"
$q4.Range("C2").Value = "华夏中证智选1000价值稳健策略ETF"
$q4.Range("C2").ClearFormats()

$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "0.36"
$q4.Range("D2").ClearFormats()

$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "96.22"
$q4.Range("E2").ClearFormats()

$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "0.88"
$q4.Range("F2").ClearFormats()

$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "0.0032"
$q4.Range("G2").ClearFormats()

$q4.Range("H2").Value = 9

# Restore the originally active sheet/tab (editing other sheets along the way
# shifts Excel's "active sheet" pointer as a side effect).
$wb.Worksheets.Item("2021-Q2").Activate()
